$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the Turkish default "Sayfa1" to "Sheet1"
$ws.Name = "Sheet1"

# Re-apply the "Normal" cell style across the used data range so Excel
# writes out explicit style records (apply* flags) for these cells
$ws.Range("A1:B4").Style = "Normal"

# Move the active selection to E22
$ws.Range("E22").Select() | Out-Null
